$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newDate = "Feb 12, 2022 (04:14:47 EST)"

# Update the Date column (F) for all data rows (2-13) to the new timestamp
for ($r = 2; $r -le 13; $r++) {
    $ws.Cells.Item($r, 6).Value = $newDate
}

# Row 7's Name cell changed from "BOMB" to "BOM"
$ws.Cells.Item(7, 2).Value = "BOM"
